$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 659.8
$ws.Range("J2").Value = 1499.5
$ws.Range("L2").Value = 1499.5
$ws.Range("N2").Value = -1725.5
$ws.Range("H15").Value = 1764.359
$ws.Range("I15").Value = 1764.359
$ws.Range("K15").Value = 5293.076999999999
$ws.Range("M15").Value = -5124.076999999999
$ws.Range("H17").Value = 1352.9412
$ws.Range("J17").Value = 1844.4445
$ws.Range("L17").Value = 5533.333500000001
$ws.Range("N17").Value = -5869.333500000001
$ws.Range("H29").Value = 117.666664
$ws.Range("I29").Value = 117.666664
$ws.Range("K29").Value = 352.999992
$ws.Range("M29").Value = -71.99999200000002
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("H38").Value = 40.6
$ws.Range("I38").Value = 40.6
$ws.Range("K38").Value = 121.8
$ws.Range("M38").Value = 250.2
$ws.Range("H58").Value = 1179.8889
$ws.Range("J58").Value = 1750
$ws.Range("L58").Value = 5250
$ws.Range("N58").Value = -5550
$ws.Range("H64").Value = 6307.143
$ws.Range("I64").Value = 5383.3335
$ws.Range("J64").Value = 7000
$ws.Range("K64").Value = 5383.3335
$ws.Range("L64").Value = 7000
$ws.Range("M64").Value = -5135.3335
$ws.Range("N64").Value = -7496
$ws.Range("H67").Value = 6307.143
$ws.Range("I67").Value = 5383.3335
$ws.Range("J67").Value = 7000
$ws.Range("K67").Value = 5383.3335
$ws.Range("L67").Value = 7000
$ws.Range("M67").Value = -4525.3335
$ws.Range("N67").Value = -8716
$ws.Range("H100").Value = 442
$ws.Range("J100").Value = 150
$ws.Range("L100").Value = 150
$ws.Range("N100").Value = -1232
$ws.Range("H116").Value = 13274.375
$ws.Range("J116").Value = 9400
$ws.Range("L116").Value = 9400
$ws.Range("N116").Value = -16284
$ws.Range("H132").Value = 3196.8462
$ws.Range("I132").Value = 1659.4375
$ws.Range("K132").Value = 4978.3125
$ws.Range("M132").Value = -2448.3125
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1200
$ws.Range("I2").Value = 1200
$ws.Range("K2").Value = 1200
$ws.Range("M2").Value = -1087
$ws.Range("H74").Value = 2313.7144
$ws.Range("I74").Value = 2073.652
$ws.Range("K74").Value = 2073.652
$ws.Range("M74").Value = -1199.652
$ws.Range("H77").Value = 2313.7144
$ws.Range("I77").Value = 2073.652
$ws.Range("K77").Value = 10368.26
$ws.Range("M77").Value = -6000.26
$ws.Range("H116").Value = 1200
$ws.Range("I116").Value = 1200
$ws.Range("K116").Value = 1200
$ws.Range("M116").Value = 1094
$ws.Range("H132").Value = 3149.3845
$ws.Range("I132").Value = 3149.3845
$ws.Range("K132").Value = 9448.1535
$ws.Range("M132").Value = -6918.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1200
$ws.Range("I3").Value = 1200
$ws.Range("K3").Value = 1200
$ws.Range("M3").Value = -1086

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7061.8
$ws.Range("I31").Value = 3273.1667
$ws.Range("J31").Value = 8685.5
$ws.Range("K31").Value = 3273.1667
$ws.Range("L31").Value = 8685.5
$ws.Range("M31").Value = -2978.1667
$ws.Range("N31").Value = -9275.5
$ws.Range("H34").Value = 7061.8
$ws.Range("I34").Value = 3273.1667
$ws.Range("J34").Value = 8685.5
$ws.Range("K34").Value = 3273.1667
$ws.Range("L34").Value = 8685.5
$ws.Range("M34").Value = -3071.1667
$ws.Range("N34").Value = -9089.5
$ws.Range("H58").Value = 3264.0667
$ws.Range("I58").Value = 1996.75
$ws.Range("K58").Value = 1996.75
$ws.Range("M58").Value = -1793.75
$ws.Range("H132").Value = 2605.8262
$ws.Range("I132").Value = 2371.95
$ws.Range("K132").Value = 7115.849999999999
$ws.Range("M132").Value = -4585.849999999999
$ws.Range("H134").Value = 4533.3335
$ws.Range("J134").Value = 5271.4287
$ws.Range("L134").Value = 15814.2861
$ws.Range("N134").Value = -20884.2861
$ws.Range("H136").Value = 3264.0667
$ws.Range("I136").Value = 1996.75
$ws.Range("K136").Value = 5990.25
$ws.Range("M136").Value = -3440.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2130.6667
$ws.Range("I34").Value = 243
$ws.Range("J34").Value = 2319.4333
$ws.Range("K34").Value = 729
$ws.Range("L34").Value = 6958.2999
$ws.Range("M34").Value = -645
$ws.Range("N34").Value = -7126.2999
$ws.Range("H38").Value = 246.75
$ws.Range("I38").Value = 264.62964
$ws.Range("K38").Value = 793.88892
$ws.Range("M38").Value = -446.88892
$ws.Range("H39").Value = 7942.4614
$ws.Range("J39").Value = 7942.4614
$ws.Range("L39").Value = 23827.3842
$ws.Range("N39").Value = -24415.3842
$ws.Range("H55").Value = 3898.2
$ws.Range("J55").Value = 4077.5789
$ws.Range("L55").Value = 12232.7367
$ws.Range("N55").Value = -12586.7367

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2291.4167
$ws.Range("I80").Value = 2637.125
$ws.Range("K80").Value = 2637.125
$ws.Range("M80").Value = -1639.125
$ws.Range("H83").Value = 2291.4167
$ws.Range("I83").Value = 2637.125
$ws.Range("K83").Value = 13185.625
$ws.Range("M83").Value = -8193.625
$ws.Range("H132").Value = 71219.60000000001
$ws.Range("I132").Value = 104150.4
$ws.Range("J132").Value = 5358
$ws.Range("K132").Value = 312451.2
$ws.Range("L132").Value = 16074
$ws.Range("M132").Value = -309921.2
$ws.Range("N132").Value = -21134
$ws.Range("H136").Value = 74775.336
$ws.Range("J136").Value = 74775.336
$ws.Range("L136").Value = 224326.008
$ws.Range("N136").Value = -229426.008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1980.25
$ws.Range("I55").Value = 1752.2858
$ws.Range("J55").Value = 2299.4
$ws.Range("K55").Value = 1752.2858
$ws.Range("L55").Value = 2299.4
$ws.Range("M55").Value = -1579.2858
$ws.Range("N55").Value = -2645.4
$ws.Range("H82").Value = 3432.9167
$ws.Range("I82").Value = 1170.7142
$ws.Range("K82").Value = 1170.7142
$ws.Range("M82").Value = -809.7141999999999
$ws.Range("H85").Value = 3432.9167
$ws.Range("I85").Value = 1170.7142
$ws.Range("K85").Value = 1170.7142
$ws.Range("M85").Value = 77.28580000000011
$ws.Range("H122").Value = 3916.6667
$ws.Range("I122").Value = 3916.6667
$ws.Range("K122").Value = 11750.0001
$ws.Range("M122").Value = -9300.000100000001
$ws.Range("H136").Value = 3999.8
$ws.Range("I136").Value = 3999.8
$ws.Range("K136").Value = 11999.4
$ws.Range("M136").Value = -9449.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10875000
$ws.Range("I5").Value = 20000000
$ws.Range("K5").Value = 20000000
$ws.Range("M5").Value = -19999888
$ws.Range("H136").Value = 3033.6128
$ws.Range("I136").Value = 2255.348
$ws.Range("K136").Value = 6766.044
$ws.Range("M136").Value = -4216.044
